# "Log In screen change" — flip the environment used by the login/asset
# screen from PROD to DEV.
#
# Settings!B3 is the single source value ("PROD/P004_NewHireCommunication").
# Settings!B20, Settings!B31 and Assets!C2:C27 all reference it (directly or
# via "Settings!B3 & "/SP_003_WorkdayDisposition""), so updating B3 alone
# ripples through every dependent formula on recalculation.
#
# The edit also leaves the Settings sheet scrolled to the top with cell B4
# selected (previously it was scrolled to row 8 with B46 selected).

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")

# Root value change: PROD -> DEV. Everything else recalculates from this.
$settings.Range("B3").Value = "DEV/P004_NewHireCommunication"

# Restore the view: scroll back to the top-left and move the selection to B4.
$settings.Activate()
$settings.Application.ActiveWindow.ScrollRow = 1
$settings.Application.ActiveWindow.ScrollColumn = 1
$settings.Range("B4").Select()
